$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim trailing whitespace from "Baz " -> "Baz"
$ws.Range("A2").Value2 = "Baz"

# Add a thin border around B2 (creates a new border + cellXf)
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("B2").Borders.Weight = 2

# Row 2's rendered height is recalculated once its formatting changes
$ws.Rows.Item(2).RowHeight = 14.9

# Move the selection from B1 to B2
$ws.Range("B2").Select()
